# Update "days since/until" counters in the F column on three sheets
# (展览, 演出, 全部类型) to reflect the new generation timestamp.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 3679
$ws1.Range("F8").Value  = 6
$ws1.Range("F9").Value  = 177
$ws1.Range("F12").Value = 1363
$ws1.Range("F14").Value = 2079

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 5

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 3679
$ws4.Range("F8").Value  = 6
$ws4.Range("F10").Value = 177
$ws4.Range("F14").Value = 5
$ws4.Range("F15").Value = 1363
$ws4.Range("F17").Value = 2079
